# Fix "New Volume" (column H) values on Sheet1.
# The previous calculation undercounted some of the per-well transfer
# volumes; this updates column H to the corrected values.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("H2").Value  = "58.2,60.0,"
$ws.Range("H3").Value  = "58.5,60.0,"
$ws.Range("H4").Value  = "50.25,"
$ws.Range("H5").Value  = "58.5,60.0,"
$ws.Range("H6").Value  = "52.5,60.0,60.0,60.0,60.0,60.0,"
$ws.Range("H7").Value  = "56.1,"
$ws.Range("H8").Value  = "44.4,60.0,"
$ws.Range("H9").Value  = "57.0,60.0,"
$ws.Range("H10").Value = "60.0,60.0,60.0,"
$ws.Range("H11").Value = "40.5,60.0,"
$ws.Range("H12").Value = "60.0,52.0,53.0,"
$ws.Range("H13").Value = "61.225,"
$ws.Range("H14").Value = "35.25,"
$ws.Range("H15").Value = "40.5,"
$ws.Range("H16").Value = "40.5,"
$ws.Range("H17").Value = "57.0,60.0,"
$ws.Range("H18").Value = "60.0,60.0,"
$ws.Range("H19").Value = "48.0,"
$ws.Range("H20").Value = "57.0,"
$ws.Range("H21").Value = "19.7999999999999,"
$ws.Range("H22").Value = "48.45,"
$ws.Range("H23").Value = "60.0,"
$ws.Range("H24").Value = "60.0,"
$ws.Range("H25").Value = "1965.05,2000.0,"
$ws.Range("H26").Value = "60.0,"
$ws.Range("H27").Value = "52.5,"
$ws.Range("H28").Value = "60.0,"
$ws.Range("H29").Value = "60.0,"
$ws.Range("H30").Value = "60.0,"
